# Rename worksheets and update the "Education[T.Unknown]" label to
# "Education[T.Unknown/Other]" on every sheet (row 5, column A).

$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ00359534",
    "summ00684686",
    "summ01056206",
    "summ01429554",
    "summ01821064",
    "summ02168926",
    "summ02512122",
    "summ02889520",
    "summ03294631"
)

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Name = $newNames[$i - 1]

    if ($ws.Range("A5").Text -eq "Education[T.Unknown]") {
        $ws.Range("A5").Value = "Education[T.Unknown/Other]"
    }
}
